$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 192,6

$data[0,0] = 45918
$data[0,1] = 0
$data[0,2] = 0
$data[0,3] = 0
$data[0,4] = 1
$data[0,5] = "18.09.20251"
$data[1,0] = 45918.01041666666
$data[1,1] = 0
$data[1,2] = 0
$data[1,3] = 0
$data[1,4] = 2
$data[1,5] = "18.09.20252"
$data[2,0] = 45918.02083333334
$data[2,1] = 0
$data[2,2] = 0
$data[2,3] = 0
$data[2,4] = 3
$data[2,5] = "18.09.20253"
$data[3,0] = 45918.03125
$data[3,1] = 0
$data[3,2] = 0
$data[3,3] = 0
$data[3,4] = 4
$data[3,5] = "18.09.20254"
$data[4,0] = 45918.04166666666
$data[4,1] = 63
$data[4,2] = -60.44
$data[4,3] = -60.44
$data[4,4] = 5
$data[4,5] = "18.09.20255"
$data[5,0] = 45918.05208333334
$data[5,1] = 28
$data[5,2] = -33.88
$data[5,3] = -33.88
$data[5,4] = 6
$data[5,5] = "18.09.20256"
$data[6,0] = 45918.0625
$data[6,1] = 46
$data[6,2] = -554.4
$data[6,3] = -554.4
$data[6,4] = 7
$data[6,5] = "18.09.20257"
$data[7,0] = 45918.07291666666
$data[7,1] = 43
$data[7,2] = -308.42
$data[7,3] = -308.42
$data[7,4] = 8
$data[7,5] = "18.09.20258"
$data[8,0] = 45918.08333333334
$data[8,1] = 49
$data[8,2] = -229.35
$data[8,3] = -229.35
$data[8,4] = 9
$data[8,5] = "18.09.20259"
$data[9,0] = 45918.09375
$data[9,1] = 12
$data[9,2] = -31.42
$data[9,3] = -31.42
$data[9,4] = 10
$data[9,5] = "18.09.202510"
$data[10,0] = 45918.10416666666
$data[10,1] = 7
$data[10,2] = 132.14
$data[10,3] = 132.14
$data[10,4] = 11
$data[10,5] = "18.09.202511"
$data[11,0] = 45918.11458333334
$data[11,1] = 14
$data[11,2] = 20.27
$data[11,3] = 20.27
$data[11,4] = 12
$data[11,5] = "18.09.202512"
$data[12,0] = 45918.125
$data[12,1] = 22
$data[12,2] = -3.67
$data[12,3] = -3.67
$data[12,4] = 13
$data[12,5] = "18.09.202513"
$data[13,0] = 45918.13541666666
$data[13,1] = 33
$data[13,2] = -14.96
$data[13,3] = -14.96
$data[13,4] = 14
$data[13,5] = "18.09.202514"
$data[14,0] = 45918.14583333334
$data[14,1] = -5
$data[14,2] = 208.73
$data[14,3] = 208.73
$data[14,4] = 15
$data[14,5] = "18.09.202515"
$data[15,0] = 45918.15625
$data[15,1] = 22
$data[15,2] = 105.39
$data[15,3] = 105.39
$data[15,4] = 16
$data[15,5] = "18.09.202516"
$data[16,0] = 45918.16666666666
$data[16,1] = 28
$data[16,2] = 149.67
$data[16,3] = 149.67
$data[16,4] = 17
$data[16,5] = "18.09.202517"
$data[17,0] = 45918.17708333334
$data[17,1] = 21
$data[17,2] = 192.44
$data[17,3] = 192.44
$data[17,4] = 18
$data[17,5] = "18.09.202518"
$data[18,0] = 45918.1875
$data[18,1] = 5
$data[18,2] = 364.7
$data[18,3] = 364.7
$data[18,4] = 19
$data[18,5] = "18.09.202519"
$data[19,0] = 45918.19791666666
$data[19,1] = 8
$data[19,2] = 265.05
$data[19,3] = 265.05
$data[19,4] = 20
$data[19,5] = "18.09.202520"
$data[20,0] = 45918.20833333334
$data[20,1] = 69
$data[20,2] = -13.43
$data[20,3] = -13.43
$data[20,4] = 21
$data[20,5] = "18.09.202521"
$data[21,0] = 45918.21875
$data[21,1] = 38
$data[21,2] = 123.56
$data[21,3] = 123.56
$data[21,4] = 22
$data[21,5] = "18.09.202522"
$data[22,0] = 45918.22916666666
$data[22,1] = 1
$data[22,2] = 166.67
$data[22,3] = 166.67
$data[22,4] = 23
$data[22,5] = "18.09.202523"
$data[23,0] = 45918.23958333334
$data[23,1] = 8
$data[23,2] = 200.98
$data[23,3] = 200.98
$data[23,4] = 24
$data[23,5] = "18.09.202524"
$data[24,0] = 45918.25
$data[24,1] = 26
$data[24,2] = -1676.86
$data[24,3] = -1676.86
$data[24,4] = 25
$data[24,5] = "18.09.202525"
$data[25,0] = 45918.26041666666
$data[25,1] = -9
$data[25,2] = 784.22
$data[25,3] = 784.22
$data[25,4] = 26
$data[25,5] = "18.09.202526"
$data[26,0] = 45918.27083333334
$data[26,1] = 2
$data[26,2] = 319.62
$data[26,3] = 319.62
$data[26,4] = 27
$data[26,5] = "18.09.202527"
$data[27,0] = 45918.28125
$data[27,1] = -3
$data[27,2] = 739.78
$data[27,3] = 739.78
$data[27,4] = 28
$data[27,5] = "18.09.202528"
$data[28,0] = 45918.29166666666
$data[28,1] = -21
$data[28,2] = 801.27
$data[28,3] = 801.27
$data[28,4] = 29
$data[28,5] = "18.09.202529"
$data[29,0] = 45918.30208333334
$data[29,1] = 38
$data[29,2] = 195.44
$data[29,3] = 195.44
$data[29,4] = 30
$data[29,5] = "18.09.202530"
$data[30,0] = 45918.3125
$data[30,1] = 67
$data[30,2] = -2056.81
$data[30,3] = -2056.81
$data[30,4] = 31
$data[30,5] = "18.09.202531"
$data[31,0] = 45918.32291666666
$data[31,1] = 93
$data[31,2] = -5999.48
$data[31,3] = -5999.48
$data[31,4] = 32
$data[31,5] = "18.09.202532"
$data[32,0] = 45918.33333333334
$data[32,1] = 16
$data[32,2] = -2347.04
$data[32,3] = -2347.04
$data[32,4] = 33
$data[32,5] = "18.09.202533"
$data[33,0] = 45918.34375
$data[33,1] = 63
$data[33,2] = -3170.71
$data[33,3] = -3170.71
$data[33,4] = 34
$data[33,5] = "18.09.202534"
$data[34,0] = 45918.35416666666
$data[34,1] = 109
$data[34,2] = -6842.96
$data[34,3] = -6842.96
$data[34,4] = 35
$data[34,5] = "18.09.202535"
$data[35,0] = 45918.36458333334
$data[35,1] = 140
$data[35,2] = -6988.34
$data[35,3] = -6988.34
$data[35,4] = 36
$data[35,5] = "18.09.202536"
$data[36,0] = 45918.375
$data[36,1] = 83
$data[36,2] = -1764.51
$data[36,3] = -1764.51
$data[36,4] = 37
$data[36,5] = "18.09.202537"
$data[37,0] = 45918.38541666666
$data[37,1] = 85
$data[37,2] = -228.77
$data[37,3] = -228.77
$data[37,4] = 38
$data[37,5] = "18.09.202538"
$data[38,0] = 45918.39583333334
$data[38,1] = 57
$data[38,2] = -470.12
$data[38,3] = -470.12
$data[38,4] = 39
$data[38,5] = "18.09.202539"
$data[39,0] = 45918.40625
$data[39,1] = 56
$data[39,2] = -1193.13
$data[39,3] = -1193.13
$data[39,4] = 40
$data[39,5] = "18.09.202540"
$data[40,0] = 45918.41666666666
$data[40,1] = 22
$data[40,2] = -151.41
$data[40,3] = -151.41
$data[40,4] = 41
$data[40,5] = "18.09.202541"
$data[41,0] = 45918.42708333334
$data[41,1] = 40
$data[41,2] = -99.58
$data[41,3] = -99.58
$data[41,4] = 42
$data[41,5] = "18.09.202542"
$data[42,0] = 45918.4375
$data[42,1] = 67
$data[42,2] = -6.81
$data[42,3] = -6.81
$data[42,4] = 43
$data[42,5] = "18.09.202543"
$data[43,0] = 45918.44791666666
$data[43,1] = 81
$data[43,2] = -3555.25
$data[43,3] = -3555.25
$data[43,4] = 44
$data[43,5] = "18.09.202544"
$data[44,0] = 45918.45833333334
$data[44,1] = 79
$data[44,2] = -2418.22
$data[44,3] = -2418.22
$data[44,4] = 45
$data[44,5] = "18.09.202545"
$data[45,0] = 45918.46875
$data[45,1] = 19
$data[45,2] = -87.31
$data[45,3] = -87.31
$data[45,4] = 46
$data[45,5] = "18.09.202546"
$data[46,0] = 45918.47916666666
$data[46,1] = -9
$data[46,2] = 450.03
$data[46,3] = 450.03
$data[46,4] = 47
$data[46,5] = "18.09.202547"
$data[47,0] = 45918.48958333334
$data[47,1] = -13
$data[47,2] = 450.03
$data[47,3] = 450.03
$data[47,4] = 48
$data[47,5] = "18.09.202548"
$data[48,0] = 45918.5
$data[48,1] = 27
$data[48,2] = 43.95
$data[48,3] = 43.95
$data[48,4] = 49
$data[48,5] = "18.09.202549"
$data[49,0] = 45918.51041666666
$data[49,1] = 20
$data[49,2] = -270.3
$data[49,3] = -270.3
$data[49,4] = 50
$data[49,5] = "18.09.202550"
$data[50,0] = 45918.52083333334
$data[50,1] = 14
$data[50,2] = -790.0599999999999
$data[50,3] = -790.0599999999999
$data[50,4] = 51
$data[50,5] = "18.09.202551"
$data[51,0] = 45918.53125
$data[51,1] = -13
$data[51,2] = 450.01
$data[51,3] = 450.01
$data[51,4] = 52
$data[51,5] = "18.09.202552"
$data[52,0] = 45918.54166666666
$data[52,1] = 32
$data[52,2] = -133.02
$data[52,3] = -133.02
$data[52,4] = 53
$data[52,5] = "18.09.202553"
$data[53,0] = 45918.55208333334
$data[53,1] = 49
$data[53,2] = 19.86
$data[53,3] = 19.86
$data[53,4] = 54
$data[53,5] = "18.09.202554"
$data[54,0] = 45918.5625
$data[54,1] = 24
$data[54,2] = -144.54
$data[54,3] = -144.54
$data[54,4] = 55
$data[54,5] = "18.09.202555"
$data[55,0] = 45918.57291666666
$data[55,1] = 42
$data[55,2] = -3036.04
$data[55,3] = -3036.04
$data[55,4] = 56
$data[55,5] = "18.09.202556"
$data[56,0] = 45918.58333333334
$data[56,1] = 91
$data[56,2] = -6983.7
$data[56,3] = -6983.7
$data[56,4] = 57
$data[56,5] = "18.09.202557"
$data[57,0] = 45918.59375
$data[57,1] = 74
$data[57,2] = -3027.57
$data[57,3] = -3027.57
$data[57,4] = 58
$data[57,5] = "18.09.202558"
$data[58,0] = 45918.60416666666
$data[58,1] = 72
$data[58,2] = -91.78
$data[58,3] = -91.78
$data[58,4] = 59
$data[58,5] = "18.09.202559"
$data[59,0] = 45918.61458333334
$data[59,1] = 109
$data[59,2] = -3477.89
$data[59,3] = -3477.89
$data[59,4] = 60
$data[59,5] = "18.09.202560"
$data[60,0] = 45918.625
$data[60,1] = 102
$data[60,2] = -6547.24
$data[60,3] = -6547.24
$data[60,4] = 61
$data[60,5] = "18.09.202561"
$data[61,0] = 45918.63541666666
$data[61,1] = 84
$data[61,2] = -4100.2
$data[61,3] = -4100.2
$data[61,4] = 62
$data[61,5] = "18.09.202562"
$data[62,0] = 45918.64583333334
$data[62,1] = 120
$data[62,2] = -6792.87
$data[62,3] = -6792.87
$data[62,4] = 63
$data[62,5] = "18.09.202563"
$data[63,0] = 45918.65625
$data[63,1] = 30
$data[63,2] = -3162.98
$data[63,3] = -3162.98
$data[63,4] = 64
$data[63,5] = "18.09.202564"
$data[64,0] = 45918.66666666666
$data[64,1] = 105
$data[64,2] = -6819.03
$data[64,3] = -6819.03
$data[64,4] = 65
$data[64,5] = "18.09.202565"
$data[65,0] = 45918.67708333334
$data[65,1] = 103
$data[65,2] = -6169.58
$data[65,3] = -6169.58
$data[65,4] = 66
$data[65,5] = "18.09.202566"
$data[66,0] = 45918.6875
$data[66,1] = 45
$data[66,2] = -583.47
$data[66,3] = -583.47
$data[66,4] = 67
$data[66,5] = "18.09.202567"
$data[67,0] = 45918.69791666666
$data[67,1] = -3
$data[67,2] = 473.11
$data[67,3] = 473.11
$data[67,4] = 68
$data[67,5] = "18.09.202568"
$data[68,0] = 45918.70833333334
$data[68,1] = 40
$data[68,2] = 236.4
$data[68,3] = 236.4
$data[68,4] = 69
$data[68,5] = "18.09.202569"
$data[69,0] = 45918.71875
$data[69,1] = 34
$data[69,2] = 390
$data[69,3] = 390
$data[69,4] = 70
$data[69,5] = "18.09.202570"
$data[70,0] = 45918.72916666666
$data[70,1] = 22
$data[70,2] = 390
$data[70,3] = 390
$data[70,4] = 71
$data[70,5] = "18.09.202571"
$data[71,0] = 45918.73958333334
$data[71,1] = 19
$data[71,2] = 400
$data[71,3] = 400
$data[71,4] = 72
$data[71,5] = "18.09.202572"
$data[72,0] = 45918.75
$data[72,1] = 89
$data[72,2] = -2111.16
$data[72,3] = -2111.16
$data[72,4] = 73
$data[72,5] = "18.09.202573"
$data[73,0] = 45918.76041666666
$data[73,1] = 75
$data[73,2] = 390
$data[73,3] = 390
$data[73,4] = 74
$data[73,5] = "18.09.202574"
$data[74,0] = 45918.77083333334
$data[74,1] = 22
$data[74,2] = 400
$data[74,3] = 400
$data[74,4] = 75
$data[74,5] = "18.09.202575"
$data[75,0] = 45918.78125
$data[75,1] = 25
$data[75,2] = 390
$data[75,3] = 390
$data[75,4] = 76
$data[75,5] = "18.09.202576"
$data[76,0] = 45918.79166666666
$data[76,1] = 48
$data[76,2] = 251.84
$data[76,3] = 251.84
$data[76,4] = 77
$data[76,5] = "18.09.202577"
$data[77,0] = 45918.80208333334
$data[77,1] = 76
$data[77,2] = 371.77
$data[77,3] = 371.77
$data[77,4] = 78
$data[77,5] = "18.09.202578"
$data[78,0] = 45918.8125
$data[78,1] = 76
$data[78,2] = 390
$data[78,3] = 390
$data[78,4] = 79
$data[78,5] = "18.09.202579"
$data[79,0] = 45918.82291666666
$data[79,1] = 94
$data[79,2] = 219.08
$data[79,3] = 219.08
$data[79,4] = 80
$data[79,5] = "18.09.202580"
$data[80,0] = 45918.83333333334
$data[80,1] = 72
$data[80,2] = -19.29
$data[80,3] = -19.29
$data[80,4] = 81
$data[80,5] = "18.09.202581"
$data[81,0] = 45918.84375
$data[81,1] = 73
$data[81,2] = -260.38
$data[81,3] = -260.38
$data[81,4] = 82
$data[81,5] = "18.09.202582"
$data[82,0] = 45918.85416666666
$data[82,1] = 60
$data[82,2] = -501.28
$data[82,3] = -501.28
$data[82,4] = 83
$data[82,5] = "18.09.202583"
$data[83,0] = 45918.86458333334
$data[83,1] = 50
$data[83,2] = -585.87
$data[83,3] = -585.87
$data[83,4] = 84
$data[83,5] = "18.09.202584"
$data[84,0] = 45918.875
$data[84,1] = 30
$data[84,2] = 111.91
$data[84,3] = 111.91
$data[84,4] = 85
$data[84,5] = "18.09.202585"
$data[85,0] = 45918.88541666666
$data[85,1] = 34
$data[85,2] = 247.55
$data[85,3] = 247.55
$data[85,4] = 86
$data[85,5] = "18.09.202586"
$data[86,0] = 45918.89583333334
$data[86,1] = 60
$data[86,2] = 139.3
$data[86,3] = 139.3
$data[86,4] = 87
$data[86,5] = "18.09.202587"
$data[87,0] = 45918.90625
$data[87,1] = 72
$data[87,2] = -2763.22
$data[87,3] = -2763.22
$data[87,4] = 88
$data[87,5] = "18.09.202588"
$data[88,0] = 45918.91666666666
$data[88,1] = 41
$data[88,2] = -937.38
$data[88,3] = -937.38
$data[88,4] = 89
$data[88,5] = "18.09.202589"
$data[89,0] = 45918.92708333334
$data[89,1] = 66
$data[89,2] = -72.56999999999999
$data[89,3] = -72.56999999999999
$data[89,4] = 90
$data[89,5] = "18.09.202590"
$data[90,0] = 45918.9375
$data[90,1] = 58
$data[90,2] = -66.31999999999999
$data[90,3] = -66.31999999999999
$data[90,4] = 91
$data[90,5] = "18.09.202591"
$data[91,0] = 45918.94791666666
$data[91,1] = 56
$data[91,2] = 184.72
$data[91,3] = 184.72
$data[91,4] = 92
$data[91,5] = "18.09.202592"
$data[92,0] = 45918.95833333334
$data[92,1] = 33
$data[92,2] = 116.65
$data[92,3] = 116.65
$data[92,4] = 93
$data[92,5] = "18.09.202593"
$data[93,0] = 45918.96875
$data[93,1] = 19
$data[93,2] = 190
$data[93,3] = 190
$data[93,4] = 94
$data[93,5] = "18.09.202594"
$data[94,0] = 45918.97916666666
$data[94,1] = 20
$data[94,2] = 178.73
$data[94,3] = 178.73
$data[94,4] = 95
$data[94,5] = "18.09.202595"
$data[95,0] = 45918.98958333334
$data[95,1] = 17
$data[95,2] = 7.07
$data[95,3] = 7.07
$data[95,4] = 96
$data[95,5] = "18.09.202596"
$data[96,0] = 45919
$data[96,1] = 95
$data[96,2] = -159.52
$data[96,3] = -159.52
$data[96,4] = 1
$data[96,5] = "19.09.20251"
$data[97,0] = 45919.01041666666
$data[97,1] = 59
$data[97,2] = 15.79
$data[97,3] = 15.79
$data[97,4] = 2
$data[97,5] = "19.09.20252"
$data[98,0] = 45919.02083333334
$data[98,1] = 40
$data[98,2] = -369.55
$data[98,3] = -369.55
$data[98,4] = 3
$data[98,5] = "19.09.20253"
$data[99,0] = 45919.03125
$data[99,1] = 54
$data[99,2] = -841.92
$data[99,3] = -841.92
$data[99,4] = 4
$data[99,5] = "19.09.20254"
$data[100,0] = 45919.04166666666
$data[100,1] = 92
$data[100,2] = -185.39
$data[100,3] = -185.39
$data[100,4] = 5
$data[100,5] = "19.09.20255"
$data[101,0] = 45919.05208333334
$data[101,1] = 38
$data[101,2] = -260.94
$data[101,3] = -260.94
$data[101,4] = 6
$data[101,5] = "19.09.20256"
$data[102,0] = 45919.0625
$data[102,1] = 57
$data[102,2] = -1302.29
$data[102,3] = -1302.29
$data[102,4] = 7
$data[102,5] = "19.09.20257"
$data[103,0] = 45919.07291666666
$data[103,1] = 61
$data[103,2] = -4694.98
$data[103,3] = -4694.98
$data[103,4] = 8
$data[103,5] = "19.09.20258"
$data[104,0] = 45919.08333333334
$data[104,1] = 64
$data[104,2] = -2128.08
$data[104,3] = -2128.08
$data[104,4] = 9
$data[104,5] = "19.09.20259"
$data[105,0] = 45919.09375
$data[105,1] = 40
$data[105,2] = -483.51
$data[105,3] = -483.51
$data[105,4] = 10
$data[105,5] = "19.09.202510"
$data[106,0] = 45919.10416666666
$data[106,1] = 43
$data[106,2] = -337.34
$data[106,3] = -337.34
$data[106,4] = 11
$data[106,5] = "19.09.202511"
$data[107,0] = 45919.11458333334
$data[107,1] = 36
$data[107,2] = -332.69
$data[107,3] = -332.69
$data[107,4] = 12
$data[107,5] = "19.09.202512"
$data[108,0] = 45919.125
$data[108,1] = 37
$data[108,2] = -429.83
$data[108,3] = -429.83
$data[108,4] = 13
$data[108,5] = "19.09.202513"
$data[109,0] = 45919.13541666666
$data[109,1] = 7
$data[109,2] = 185.66
$data[109,3] = 185.66
$data[109,4] = 14
$data[109,5] = "19.09.202514"
$data[110,0] = 45919.14583333334
$data[110,1] = 13
$data[110,2] = 132.18
$data[110,3] = 132.18
$data[110,4] = 15
$data[110,5] = "19.09.202515"
$data[111,0] = 45919.15625
$data[111,1] = 12
$data[111,2] = 143.73
$data[111,3] = 143.73
$data[111,4] = 16
$data[111,5] = "19.09.202516"
$data[112,0] = 45919.16666666666
$data[112,1] = 34
$data[112,2] = -309.92
$data[112,3] = -309.92
$data[112,4] = 17
$data[112,5] = "19.09.202517"
$data[113,0] = 45919.17708333334
$data[113,1] = 19
$data[113,2] = -13.97
$data[113,3] = -13.97
$data[113,4] = 18
$data[113,5] = "19.09.202518"
$data[114,0] = 45919.1875
$data[114,1] = 40
$data[114,2] = -149.43
$data[114,3] = -149.43
$data[114,4] = 19
$data[114,5] = "19.09.202519"
$data[115,0] = 45919.19791666666
$data[115,1] = 36
$data[115,2] = 24.66
$data[115,3] = 24.66
$data[115,4] = 20
$data[115,5] = "19.09.202520"
$data[116,0] = 45919.20833333334
$data[116,1] = 60
$data[116,2] = -2907.94
$data[116,3] = -2907.94
$data[116,4] = 21
$data[116,5] = "19.09.202521"
$data[117,0] = 45919.21875
$data[117,1] = 30
$data[117,2] = -201.99
$data[117,3] = -201.99
$data[117,4] = 22
$data[117,5] = "19.09.202522"
$data[118,0] = 45919.22916666666
$data[118,1] = 37
$data[118,2] = 194.11
$data[118,3] = 194.11
$data[118,4] = 23
$data[118,5] = "19.09.202523"
$data[119,0] = 45919.23958333334
$data[119,1] = 24
$data[119,2] = -199
$data[119,3] = -199
$data[119,4] = 24
$data[119,5] = "19.09.202524"
$data[120,0] = 45919.25
$data[120,1] = 25
$data[120,2] = 202.38
$data[120,3] = 202.38
$data[120,4] = 25
$data[120,5] = "19.09.202525"
$data[121,0] = 45919.26041666666
$data[121,1] = 17
$data[121,2] = 232.88
$data[121,3] = 232.88
$data[121,4] = 26
$data[121,5] = "19.09.202526"
$data[122,0] = 45919.27083333334
$data[122,1] = 32
$data[122,2] = 358.78
$data[122,3] = 358.78
$data[122,4] = 27
$data[122,5] = "19.09.202527"
$data[123,0] = 45919.28125
$data[123,1] = 18
$data[123,2] = 389.09
$data[123,3] = 389.09
$data[123,4] = 28
$data[123,5] = "19.09.202528"
$data[124,0] = 45919.29166666666
$data[124,1] = -44
$data[124,2] = 1080.2
$data[124,3] = 1080.2
$data[124,4] = 29
$data[124,5] = "19.09.202529"
$data[125,0] = 45919.30208333334
$data[125,1] = -60
$data[125,2] = 811.9
$data[125,3] = 811.9
$data[125,4] = 30
$data[125,5] = "19.09.202530"
$data[126,0] = 45919.3125
$data[126,1] = -59
$data[126,2] = 980.5599999999999
$data[126,3] = 980.5599999999999
$data[126,4] = 31
$data[126,5] = "19.09.202531"
$data[127,0] = 45919.32291666666
$data[127,1] = -40
$data[127,2] = 799.97
$data[127,3] = 799.97
$data[127,4] = 32
$data[127,5] = "19.09.202532"
$data[128,0] = 45919.33333333334
$data[128,1] = -101
$data[128,2] = 6163.32
$data[128,3] = 6163.32
$data[128,4] = 33
$data[128,5] = "19.09.202533"
$data[129,0] = 45919.34375
$data[129,1] = -133
$data[129,2] = 3133.2
$data[129,3] = 3133.2
$data[129,4] = 34
$data[129,5] = "19.09.202534"
$data[130,0] = 45919.35416666666
$data[130,1] = -142
$data[130,2] = 3279.61
$data[130,3] = 3279.61
$data[130,4] = 35
$data[130,5] = "19.09.202535"
$data[131,0] = 45919.36458333334
$data[131,1] = -63
$data[131,2] = 791.6799999999999
$data[131,3] = 791.6799999999999
$data[131,4] = 36
$data[131,5] = "19.09.202536"
$data[132,0] = 45919.375
$data[132,1] = -97
$data[132,2] = 5678.21
$data[132,3] = 5678.21
$data[132,4] = 37
$data[132,5] = "19.09.202537"
$data[133,0] = 45919.38541666666
$data[133,1] = -41
$data[133,2] = 995.49
$data[133,3] = 995.49
$data[133,4] = 38
$data[133,5] = "19.09.202538"
$data[134,0] = 45919.39583333334
$data[134,1] = -8
$data[134,2] = 628
$data[134,3] = 628
$data[134,4] = 39
$data[134,5] = "19.09.202539"
$data[135,0] = 45919.40625
$data[135,1] = 0
$data[135,2] = 0
$data[135,3] = 0
$data[135,4] = 40
$data[135,5] = "19.09.202540"
$data[136,0] = 45919.41666666666
$data[136,1] = 0
$data[136,2] = 0
$data[136,3] = 0
$data[136,4] = 41
$data[136,5] = "19.09.202541"
$data[137,0] = 45919.42708333334
$data[137,1] = 0
$data[137,2] = 0
$data[137,3] = 0
$data[137,4] = 42
$data[137,5] = "19.09.202542"
$data[138,0] = 45919.4375
$data[138,1] = 0
$data[138,2] = 0
$data[138,3] = 0
$data[138,4] = 43
$data[138,5] = "19.09.202543"
$data[139,0] = 45919.44791666666
$data[139,1] = 0
$data[139,2] = 0
$data[139,3] = 0
$data[139,4] = 44
$data[139,5] = "19.09.202544"
$data[140,0] = 45919.45833333334
$data[140,1] = 0
$data[140,2] = 0
$data[140,3] = 0
$data[140,4] = 45
$data[140,5] = "19.09.202545"
$data[141,0] = 45919.46875
$data[141,1] = 0
$data[141,2] = 0
$data[141,3] = 0
$data[141,4] = 46
$data[141,5] = "19.09.202546"
$data[142,0] = 45919.47916666666
$data[142,1] = 0
$data[142,2] = 0
$data[142,3] = 0
$data[142,4] = 47
$data[142,5] = "19.09.202547"
$data[143,0] = 45919.48958333334
$data[143,1] = 0
$data[143,2] = 0
$data[143,3] = 0
$data[143,4] = 48
$data[143,5] = "19.09.202548"
$data[144,0] = 45919.5
$data[144,1] = 0
$data[144,2] = 0
$data[144,3] = 0
$data[144,4] = 49
$data[144,5] = "19.09.202549"
$data[145,0] = 45919.51041666666
$data[145,1] = 0
$data[145,2] = 0
$data[145,3] = 0
$data[145,4] = 50
$data[145,5] = "19.09.202550"
$data[146,0] = 45919.52083333334
$data[146,1] = 0
$data[146,2] = 0
$data[146,3] = 0
$data[146,4] = 51
$data[146,5] = "19.09.202551"
$data[147,0] = 45919.53125
$data[147,1] = 0
$data[147,2] = 0
$data[147,3] = 0
$data[147,4] = 52
$data[147,5] = "19.09.202552"
$data[148,0] = 45919.54166666666
$data[148,1] = 0
$data[148,2] = 0
$data[148,3] = 0
$data[148,4] = 53
$data[148,5] = "19.09.202553"
$data[149,0] = 45919.55208333334
$data[149,1] = 0
$data[149,2] = 0
$data[149,3] = 0
$data[149,4] = 54
$data[149,5] = "19.09.202554"
$data[150,0] = 45919.5625
$data[150,1] = 0
$data[150,2] = 0
$data[150,3] = 0
$data[150,4] = 55
$data[150,5] = "19.09.202555"
$data[151,0] = 45919.57291666666
$data[151,1] = 0
$data[151,2] = 0
$data[151,3] = 0
$data[151,4] = 56
$data[151,5] = "19.09.202556"
$data[152,0] = 45919.58333333334
$data[152,1] = 0
$data[152,2] = 0
$data[152,3] = 0
$data[152,4] = 57
$data[152,5] = "19.09.202557"
$data[153,0] = 45919.59375
$data[153,1] = 0
$data[153,2] = 0
$data[153,3] = 0
$data[153,4] = 58
$data[153,5] = "19.09.202558"
$data[154,0] = 45919.60416666666
$data[154,1] = 0
$data[154,2] = 0
$data[154,3] = 0
$data[154,4] = 59
$data[154,5] = "19.09.202559"
$data[155,0] = 45919.61458333334
$data[155,1] = 0
$data[155,2] = 0
$data[155,3] = 0
$data[155,4] = 60
$data[155,5] = "19.09.202560"
$data[156,0] = 45919.625
$data[156,1] = 0
$data[156,2] = 0
$data[156,3] = 0
$data[156,4] = 61
$data[156,5] = "19.09.202561"
$data[157,0] = 45919.63541666666
$data[157,1] = 0
$data[157,2] = 0
$data[157,3] = 0
$data[157,4] = 62
$data[157,5] = "19.09.202562"
$data[158,0] = 45919.64583333334
$data[158,1] = 0
$data[158,2] = 0
$data[158,3] = 0
$data[158,4] = 63
$data[158,5] = "19.09.202563"
$data[159,0] = 45919.65625
$data[159,1] = 0
$data[159,2] = 0
$data[159,3] = 0
$data[159,4] = 64
$data[159,5] = "19.09.202564"
$data[160,0] = 45919.66666666666
$data[160,1] = 0
$data[160,2] = 0
$data[160,3] = 0
$data[160,4] = 65
$data[160,5] = "19.09.202565"
$data[161,0] = 45919.67708333334
$data[161,1] = 0
$data[161,2] = 0
$data[161,3] = 0
$data[161,4] = 66
$data[161,5] = "19.09.202566"
$data[162,0] = 45919.6875
$data[162,1] = 0
$data[162,2] = 0
$data[162,3] = 0
$data[162,4] = 67
$data[162,5] = "19.09.202567"
$data[163,0] = 45919.69791666666
$data[163,1] = 0
$data[163,2] = 0
$data[163,3] = 0
$data[163,4] = 68
$data[163,5] = "19.09.202568"
$data[164,0] = 45919.70833333334
$data[164,1] = 0
$data[164,2] = 0
$data[164,3] = 0
$data[164,4] = 69
$data[164,5] = "19.09.202569"
$data[165,0] = 45919.71875
$data[165,1] = 0
$data[165,2] = 0
$data[165,3] = 0
$data[165,4] = 70
$data[165,5] = "19.09.202570"
$data[166,0] = 45919.72916666666
$data[166,1] = 0
$data[166,2] = 0
$data[166,3] = 0
$data[166,4] = 71
$data[166,5] = "19.09.202571"
$data[167,0] = 45919.73958333334
$data[167,1] = 0
$data[167,2] = 0
$data[167,3] = 0
$data[167,4] = 72
$data[167,5] = "19.09.202572"
$data[168,0] = 45919.75
$data[168,1] = 0
$data[168,2] = 0
$data[168,3] = 0
$data[168,4] = 73
$data[168,5] = "19.09.202573"
$data[169,0] = 45919.76041666666
$data[169,1] = 0
$data[169,2] = 0
$data[169,3] = 0
$data[169,4] = 74
$data[169,5] = "19.09.202574"
$data[170,0] = 45919.77083333334
$data[170,1] = 0
$data[170,2] = 0
$data[170,3] = 0
$data[170,4] = 75
$data[170,5] = "19.09.202575"
$data[171,0] = 45919.78125
$data[171,1] = 0
$data[171,2] = 0
$data[171,3] = 0
$data[171,4] = 76
$data[171,5] = "19.09.202576"
$data[172,0] = 45919.79166666666
$data[172,1] = 0
$data[172,2] = 0
$data[172,3] = 0
$data[172,4] = 77
$data[172,5] = "19.09.202577"
$data[173,0] = 45919.80208333334
$data[173,1] = 0
$data[173,2] = 0
$data[173,3] = 0
$data[173,4] = 78
$data[173,5] = "19.09.202578"
$data[174,0] = 45919.8125
$data[174,1] = 0
$data[174,2] = 0
$data[174,3] = 0
$data[174,4] = 79
$data[174,5] = "19.09.202579"
$data[175,0] = 45919.82291666666
$data[175,1] = 0
$data[175,2] = 0
$data[175,3] = 0
$data[175,4] = 80
$data[175,5] = "19.09.202580"
$data[176,0] = 45919.83333333334
$data[176,1] = 0
$data[176,2] = 0
$data[176,3] = 0
$data[176,4] = 81
$data[176,5] = "19.09.202581"
$data[177,0] = 45919.84375
$data[177,1] = 0
$data[177,2] = 0
$data[177,3] = 0
$data[177,4] = 82
$data[177,5] = "19.09.202582"
$data[178,0] = 45919.85416666666
$data[178,1] = 0
$data[178,2] = 0
$data[178,3] = 0
$data[178,4] = 83
$data[178,5] = "19.09.202583"
$data[179,0] = 45919.86458333334
$data[179,1] = 0
$data[179,2] = 0
$data[179,3] = 0
$data[179,4] = 84
$data[179,5] = "19.09.202584"
$data[180,0] = 45919.875
$data[180,1] = 0
$data[180,2] = 0
$data[180,3] = 0
$data[180,4] = 85
$data[180,5] = "19.09.202585"
$data[181,0] = 45919.88541666666
$data[181,1] = 0
$data[181,2] = 0
$data[181,3] = 0
$data[181,4] = 86
$data[181,5] = "19.09.202586"
$data[182,0] = 45919.89583333334
$data[182,1] = 0
$data[182,2] = 0
$data[182,3] = 0
$data[182,4] = 87
$data[182,5] = "19.09.202587"
$data[183,0] = 45919.90625
$data[183,1] = 0
$data[183,2] = 0
$data[183,3] = 0
$data[183,4] = 88
$data[183,5] = "19.09.202588"
$data[184,0] = 45919.91666666666
$data[184,1] = 0
$data[184,2] = 0
$data[184,3] = 0
$data[184,4] = 89
$data[184,5] = "19.09.202589"
$data[185,0] = 45919.92708333334
$data[185,1] = 0
$data[185,2] = 0
$data[185,3] = 0
$data[185,4] = 90
$data[185,5] = "19.09.202590"
$data[186,0] = 45919.9375
$data[186,1] = 0
$data[186,2] = 0
$data[186,3] = 0
$data[186,4] = 91
$data[186,5] = "19.09.202591"
$data[187,0] = 45919.94791666666
$data[187,1] = 0
$data[187,2] = 0
$data[187,3] = 0
$data[187,4] = 92
$data[187,5] = "19.09.202592"
$data[188,0] = 45919.95833333334
$data[188,1] = 0
$data[188,2] = 0
$data[188,3] = 0
$data[188,4] = 93
$data[188,5] = "19.09.202593"
$data[189,0] = 45919.96875
$data[189,1] = 0
$data[189,2] = 0
$data[189,3] = 0
$data[189,4] = 94
$data[189,5] = "19.09.202594"
$data[190,0] = 45919.97916666666
$data[190,1] = 0
$data[190,2] = 0
$data[190,3] = 0
$data[190,4] = 95
$data[190,5] = "19.09.202595"
$data[191,0] = 45919.98958333334
$data[191,1] = 0
$data[191,2] = 0
$data[191,3] = 0
$data[191,4] = 96
$data[191,5] = "19.09.202596"

$ws.Range("A2:F193").Value = $data
